# conductivity_field.xlsx maintenance pass: fill in two logger serial
# numbers that had come in since the last save, append the newest
# reading (2020-11-17) to every site's data log, and leave the
# selection/active-sheet state where the author was last looking
# (PBSF, reviewing the new row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "loggers" sheet: fill in two previously-blank logger serial
#    numbers and clear their special ("Good"/green) cell style, which
#    was only ever a placeholder highlight for "still missing".
# ---------------------------------------------------------------------
$loggers = $wb.Worksheets.Item("loggers")

$loggers.Range("A10").Style = "Normal"
$loggers.Range("A10").Value = 20882401

$loggers.Range("A17").Style = "Normal"
$loggers.Range("A17").Value = 20882406

# ---------------------------------------------------------------------
# 2. Append one new logger reading to each site's data sheet.
#    Each row: date-serial (col A, formatted like the row above it),
#    conductivity (col B), temperature (col C).
# ---------------------------------------------------------------------

# WIC
$ws = $wb.Worksheets.Item("WIC")
$ws.Range("A12").Value = 44152.423611111109
$ws.Range("A11").Copy() | Out-Null
$null = $ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B12").Value = 480.5
$ws.Range("C12").Value = 2.6

# YS
$ws = $wb.Worksheets.Item("YS")
$ws.Range("A27").Value = 44152.438888888886
$ws.Range("A26").Copy() | Out-Null
$null = $ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B27").Value = 347.4
$ws.Range("C27").Value = 4.7

# SW - row 26 already exists (G/H/I cells present); fill in A/B/C.
$ws = $wb.Worksheets.Item("SW")
$ws.Range("A26").Value = 44152.459027777775
$ws.Range("A25").Copy() | Out-Null
$null = $ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B26").Value = 682.5
$ws.Range("C26").Value = 3.9

# YI
$ws = $wb.Worksheets.Item("YI")
$ws.Range("A25").Value = 44152.475694444445
$ws.Range("A24").Copy() | Out-Null
$null = $ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B25").Value = 343.1
$ws.Range("C25").Value = 5.9

# YN
$ws = $wb.Worksheets.Item("YN")
$ws.Range("A25").Value = 44152.529166666667
$ws.Range("A24").Copy() | Out-Null
$null = $ws.Range("A25").PasteSpecial(-4122)
$ws.Range("B25").Value = 378
$ws.Range("C25").Value = 2.8

# 6MC
$ws = $wb.Worksheets.Item("6MC")
$ws.Range("A26").Value = 44152.541666666664
$ws.Range("A25").Copy() | Out-Null
$null = $ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B26").Value = 433
$ws.Range("C26").Value = 3.5

# DC
$ws = $wb.Worksheets.Item("DC")
$ws.Range("A26").Value = 44152.554861111108
$ws.Range("A25").Copy() | Out-Null
$null = $ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B26").Value = 510
$ws.Range("C26").Value = 4.4000000000000004

# PBMS
$ws = $wb.Worksheets.Item("PBMS")
$ws.Range("A27").Value = 44152.576388888891
$ws.Range("A26").Copy() | Out-Null
$null = $ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B27").Value = 517.4
$ws.Range("C27").Value = 3

# PBSF
$ws = $wb.Worksheets.Item("PBSF")
$ws.Range("A28").Value = 44152.589583333334
$ws.Range("A27").Copy() | Out-Null
$null = $ws.Range("A28").PasteSpecial(-4122)
$ws.Range("B28").Value = 230.2
$ws.Range("C28").Value = 4.8

# ---------------------------------------------------------------------
# 3. Restore each sheet's on-screen selection/scroll position, and
#    leave "loggers"'s frozen-pane scrolled back to the left (it was
#    scrolled right while the serial numbers above were being typed).
# ---------------------------------------------------------------------

$loggers.Range("C32").Select() | Out-Null

$wb.Worksheets.Item("WIC").Range("D24").Select() | Out-Null
$wb.Worksheets.Item("YS").Range("C37").Select() | Out-Null
$wb.Worksheets.Item("SW").Range("C26").Select() | Out-Null
$wb.Worksheets.Item("YI").Range("D34").Select() | Out-Null
$wb.Worksheets.Item("YN").Range("C24").Select() | Out-Null
$wb.Worksheets.Item("6MC").Range("C26").Select() | Out-Null
$wb.Worksheets.Item("DC").Range("D26").Select() | Out-Null
$wb.Worksheets.Item("PBMS").Range("C28").Select() | Out-Null

# PBSF ends up the active tab/sheet, selection on the just-typed row.
$ws = $wb.Worksheets.Item("PBSF")
$ws.Range("B18").Select() | Out-Null
